$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.16380786895752
$ws.Range("B1").Value = 2.369027853012085
$ws.Range("D1").Value = 2.394299030303955
$ws.Range("E1").Value = 1.21688437461853
